# Added filtering options for the Component Analysis
# Removes the "extra" forecast-horizon columns (Q5-Q9, i.e. columns G-K) from rows
# where those horizons should not be reported, keeping only a staircase window of
# columns per row, matching the new filtering logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> last column letter that should remain populated.
# Any populated cell in a column AFTER this one (up to K) gets cleared for that row.
$lastColByRow = @{
    2  = "F";  3  = "H";  4  = "F";  5  = "H";  6  = "F";  7  = "H";  8  = "F";
    9  = "H"; 10  = "F"; 11  = "H"; 12  = "F"; 13  = "H"; 14  = "F"; 15  = "H";
    16 = "F"; 17  = "H"; 18  = "J"; 19  = "H"; 20  = "J"; 21  = "H"; 22  = "J";
    23 = "I"; 24  = "H"; 25  = "K"; 26  = "J"; 27  = "I"; 28  = "H"; 29  = "K";
    30 = "J"; 31  = "I"; 32  = "H"; 33  = "K"; 34  = "J"; 35  = "I"; 36  = "H";
    37 = "K"; 38  = "J"; 39  = "I"; 40  = "H"; 41  = "K"; 42  = "J"; 43  = "I";
    44 = "H"
}

# Column letter -> index (G=7 .. K=11)
$colIndex = @{ "A"=1; "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11 }

foreach ($row in $lastColByRow.Keys) {
    $lastCol = $lastColByRow[$row]
    $startIdx = $colIndex[$lastCol] + 1
    for ($c = $startIdx; $c -le 11; $c++) {
        $ws.Cells.Item($row, $c).ClearContents()
    }
}
